# The underlying dataset rows (2..19 on "Artfynd") were re-fetched from the
# source system and came back in a different (non-deterministic) order.
# Row 18 happens to stay put; every other row's whole record relocates to a
# different row number, verbatim (no field-level edits beyond the move).
#
# before-row -> after-row
$rowMap = @{
    2  = 4
    3  = 19
    4  = 12
    5  = 7
    6  = 2
    7  = 17
    8  = 3
    9  = 11
    10 = 14
    11 = 6
    12 = 8
    13 = 9
    14 = 5
    15 = 16
    16 = 10
    17 = 15
    18 = 18
    19 = 13
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstCol = 1   # A
$lastCol  = 51  # AY

# --- Pass 1: snapshot every source row into memory before writing anything,
# since the permutation has cycles and an in-place write would clobber a
# not-yet-read row.
$snapshots = @{}
foreach ($srcRow in $rowMap.Keys) {
    $rng = $ws.Range($ws.Cells.Item($srcRow, $firstCol), $ws.Cells.Item($srcRow, $lastCol))
    $snapshots[$srcRow] = $rng.Value2
}

# --- Pass 2: write each snapshot into its destination row. Numeric/boolean
# cells go back as-is; string cells get a leading apostrophe so Excel's COM
# layer can't "helpfully" reinterpret numeric-looking or date/time-looking
# text (e.g. "30", "2023-09-04", "14:38") as a real number/date, which would
# silently change the cell's type and stored value. The apostrophe itself
# leaves no literal mark in the stored text, but it does set a transient
# "quote prefix" cell style, so we reset each written row's Style back to
# Normal afterwards to avoid introducing a spurious formatting diff.
foreach ($srcRow in $rowMap.Keys) {
    $dstRow = $rowMap[$srcRow]
    $data = $snapshots[$srcRow]

    $outArr = New-Object 'object[,]' 1, $lastCol
    for ($c = 1; $c -le $lastCol; $c++) {
        $val = $data[1, $c]
        if ($val -is [string]) {
            $outArr[0, $c - 1] = "'" + $val
        } else {
            $outArr[0, $c - 1] = $val
        }
    }

    $dstRng = $ws.Range($ws.Cells.Item($dstRow, $firstCol), $ws.Cells.Item($dstRow, $lastCol))
    $dstRng.Value2 = $outArr
    $dstRng.Style = "Normal"
}
